$wb = $excel.ActiveWorkbook

# Sheet "2025" (sheet1.xml)
$ws = $wb.Worksheets.Item("2025")
$ws.Range("A2").Value = 0
$ws.Range("B2").Value = 4357.823875228254
$ws.Range("E2").Value = 288299.5429482079
$ws.Range("G2").Value = 80959.25712661834
$ws.Range("I2").Value = 159992.7815958817
$ws.Range("L2").Value = 487002.2979492
$ws.Range("M2").Value = 112861.3269883
$ws.Range("N2").Value = 71709.53159849434
$ws.Range("O2").Value = 67449.74959589262

# Sheet "2030" (sheet2.xml)
$ws = $wb.Worksheets.Item("2030")
$ws.Range("A2").Value = 4128.642372630347
$ws.Range("B2").Value = 37397.57803605858
$ws.Range("E2").Value = 157057.2103891942
$ws.Range("I2").Value = 208211.38052173
$ws.Range("L2").Value = 0
$ws.Range("M2").Value = 61591.76966197747
$ws.Range("N2").Value = 18360.80710027935
$ws.Range("O2").Value = 10996.14522950405

# Sheet "2035" (sheet3.xml)
$ws = $wb.Worksheets.Item("2035")
$ws.Range("A2").Value = 20927.37153491941
$ws.Range("B2").Value = 15929.09817745369
$ws.Range("E2").Value = 117187.0091721428
$ws.Range("I2").Value = 167438.2047553067
$ws.Range("L2").Value = 0
$ws.Range("M2").Value = 64323.74191122008
$ws.Range("N2").Value = 43698.25102198371
$ws.Range("O2").Value = 50960.0028627053
